$d = $word.ActiveDocument

function Break-At($findText, $replaceText) {
    $range = $d.Content
    $range.Find.ClearFormatting()
    $range.Find.Replacement.ClearFormatting()
    # MatchCase = $true (2nd arg) so runs differing only by case are not conflated.
    $ok = $range.Find.Execute($findText, $true, $false, $false, $false, $false, $true, 1, $false, $replaceText, 2)
    if (-not $ok) {
        throw "Find.Execute failed for: $findText"
    }
}

# --- Portuguese "Programa" paragraph ---
Break-At "sistemas).Estrutura" "sistemas).^lEstrutura"
Break-At "Periódica.A Ligação" "Periódica.^lA Ligação"
Break-At "ligação. Natureza dos Compostos" "ligação. ^lNatureza dos Compostos"
Break-At "intermoleculares.Reações Químicas" "intermoleculares.^lReações Químicas"
Break-At "Redução).Gases: Variáveis" "Redução).^lGases: Variáveis"
Break-At "Avogadro.Soluções:" "Avogadro.^lSoluções:"
Break-At "temperatura.Estequiometria" "temperatura.^lEstequiometria"

# --- English "Programa" paragraph ---
Break-At "systems).Atomic" "systems).^lAtomic"
Break-At "Table.The Chemical" "Table.^lThe Chemical"
Break-At "bonding. Nature of the Compounds" "bonding. ^lNature of the Compounds"
Break-At "forces.Chemical Reactions" "forces.^lChemical Reactions"
Break-At "Reduction).Gases: State" "Reduction).^lGases: State"
Break-At "Principle.Solutions:" "Principle.^lSolutions:"
Break-At "temperature.Stoichiometric" "temperature.^lStoichiometric"

# --- Bibliography paragraph ---
Break-At "2005-2007ATKINS" "2005-2007^lATKINS"
Break-At "2006BRADY" "2006^lBRADY"
Break-At "1981CHANG" "1981^lCHANG"
Break-At "2010.RUSSEL" "2010.^lRUSSEL"

Write-Host "Done"
